$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 gets an explicit custom row height (19.5pt)
$ws.Rows.Item(1).RowHeight = 19.5

# Row 3 (header top border row): extend the bordered range into the new
# column Q by copying the formatting from the existing O3/P3 cell.
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 4 (year headers): add 2020 under a new column Q, formatted like P4.
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("Q4").Value = 2020

# Row 5 data value for 2020
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("Q5").Value = 25.6

# Row 6 data value for 2020
$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("Q6").Value = 13.073527219449954

# Row 7 data value for 2020
$ws.Range("P7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("Q7").Value = 21.941290626870046

# Row 8 data value for 2020
$ws.Range("P8").Copy()
$ws.Range("Q8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("Q8").Value = 196.6
